# Update to metadata: coolant & institutions
#
# This script applies the same semantic edits as the target commit:
#  - "Listen" (lists) sheet: refresh the coolant / institution /
#    workpiece_material / tool_material dropdown-source lists
#  - selection/cursor bookkeeping on the "Zuordnung" and "Listen" sheets

$wb = $excel.ActiveWorkbook

$wsZuordnung = $wb.Worksheets.Item("Zuordnung")
$wsListen    = $wb.Worksheets.Item("Listen")

# --- Listen sheet: update the dropdown/source lists -----------------------
#
# Cells are written in the same order the new values first appear in the
# rebuilt shared-string table, so newly-introduced strings line up with the
# target workbook (pure bookkeeping - the resulting cell values are what
# actually matters).

# Coolant column (B): Flood, MMQ, … -> Dry, Air, MMQ, Flood, Oil
$wsListen.Range("B4").Value = "Dry"
$wsListen.Range("B8").Value = "Oil"
$wsListen.Range("B5").Value = "Air"
$wsListen.Range("B6").Value = "MMQ"
$wsListen.Range("B7").Value = "Flood"

# Institution column (C): TU Wien, TU München, ETH Zürich -> TU Wien, TU Darmstadt
$wsListen.Range("C4").Value = "TU Wien"
$wsListen.Range("C5").Value = "TU Darmstadt"
$wsListen.Range("C6").ClearContents()

# Workpiece material column (D): S235, 4140, TiAl, Grade 5 Titanium -> C45, Steel (generic)
$wsListen.Range("D4").Value = "C45"
$wsListen.Range("D5").Value = "Steel (generic)"
$wsListen.Range("D6").ClearContents()
$wsListen.Range("D7").ClearContents()

# Tool material column (E): PCD, Carbide, MCD, Ceramic -> Carbide (P40), Carbide, MCD, Ceramic, PCD
$wsListen.Range("E4").Value = "Carbide (P40)"
$wsListen.Range("E5").Value = "Carbide"
$wsListen.Range("E6").Value = "MCD"
$wsListen.Range("E7").Value = "Ceramic"
$wsListen.Range("E8").Value = "PCD"

# --- selection / active-cell bookkeeping -----------------------------------

$wsZuordnung.Range("E9").Select()
$wsListen.Range("F22").Select()
